$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "reviews_count" column (column E), shifting columns F:K left by one
$ws.Range("E1:E1").EntireColumn.Delete()
